# FOHM covid-19 dashboard data refresh: 2020-06-06 -> 2020-06-07
# Mirrors the upstream diff: revises the last couple of daily rows on the
# three "per dag" sheets, appends a new day (2020-06-07 / serial 43989),
# refreshes the four "Totalt" aggregate sheets, moves the "Uppgift saknas"
# bucket row to the bottom of the "Antal avlidna per dag" sheet, and
# renames the FOHM snapshot tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Antal per dag region"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Revise row 124 (2020-06-05 / 43987) and row 125 (2020-06-06 / 43988)
# with the latest backfilled figures.
$row124 = @{1=43987;2=1140;3=10;4=20;5=1;6=43;7=23;8=13;9=42;10=14;11=19;12=16;13=38;14=236;15=3;16=34;17=12;18=5;19=22;20=19;21=483;22=46;23=41}
foreach ($col in $row124.Keys) {
    $ws1.Cells.Item(124, $col).Value = $row124[$col]
}

$row125 = @{1=43988;2=752;3=5;4=16;5=1;6=32;7=19;8=22;9=40;10=1;11=0;12=7;13=23;14=188;15=1;16=42;17=24;18=8;19=7;20=0;21=257;22=15;23=44}
foreach ($col in $row125.Keys) {
    $ws1.Cells.Item(125, $col).Value = $row125[$col]
}

# Append new row 126 (2020-06-07 / 43989); clone formatting from row 125.
$ws1.Range("A125:W125").Copy()
$ws1.Range("A126:W126").PasteSpecial(-4122)
$row126 = @{1=43989;2=81;3=0;4=1;5=0;6=2;7=6;8=0;9=0;10=1;11=0;12=0;13=2;14=12;15=0;16=9;17=2;18=1;19=0;20=0;21=44;22=0;23=1}
foreach ($col in $row126.Keys) {
    $ws1.Cells.Item(126, $col).Value = $row126[$col]
}

# Column A width tweak (pixel-quantised by the host, closest achievable).
$ws1.Columns.Item(1).ColumnWidth = 13.7

# ---------------------------------------------------------------------
# Sheet 2: "Antal avlidna per dag"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(81, 2).Value = 34
$ws2.Cells.Item(87, 2).Value = 11

# Row 89 used to hold the "Uppgift saknas" bucket; it now becomes the
# 2020-06-06 (43988) daily row, and the bucket moves to a new row 90.
$ws2.Range("A88:B88").Copy()
$ws2.Range("A90:B90").PasteSpecial(-4122)
$ws2.Cells.Item(89, 1).Value = 43988
$ws2.Cells.Item(89, 2).Value = 1
$ws2.Cells.Item(90, 1).Value = "Uppgift saknas"
$ws2.Cells.Item(90, 2).Value = 10

# ---------------------------------------------------------------------
# Sheet 3: "Antal intensivvårdade per dag"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(93, 2).Value = 8

$ws3.Range("A93:B93").Copy()
$ws3.Range("A94:B95").PasteSpecial(-4122)
$ws3.Cells.Item(94, 1).Value = 43988
$ws3.Cells.Item(94, 2).Value = 6
$ws3.Cells.Item(95, 1).Value = 43989
$ws3.Cells.Item(95, 2).Value = 1

# ---------------------------------------------------------------------
# Sheet 4: "Totalt antal per region"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2, 2).Value = 279
$ws4.Cells.Item(2, 3).Value = 174.80546569824219

$ws4.Cells.Item(3, 2).Value = 1314
$ws4.Cells.Item(3, 3).Value = 456.30386352539063

$ws4.Cells.Item(4, 2).Value = 90
$ws4.Cells.Item(4, 3).Value = 150.78912353515625

$ws4.Cells.Item(5, 2).Value = 1503
$ws4.Cells.Item(5, 3).Value = 522.99725341796875
$ws4.Cells.Item(5, 4).Value = 58

$ws4.Cells.Item(6, 2).Value = 913
$ws4.Cells.Item(6, 3).Value = 273.47775268554688

$ws4.Cells.Item(7, 2).Value = 803
$ws4.Cells.Item(7, 3).Value = 613.867431640625
$ws4.Cells.Item(7, 5).Value = 34

$ws4.Cells.Item(8, 2).Value = 1750
$ws4.Cells.Item(8, 3).Value = 481.2994384765625
$ws4.Cells.Item(8, 4).Value = 78

$ws4.Cells.Item(9, 2).Value = 406
$ws4.Cells.Item(9, 3).Value = 165.41316223144531

$ws4.Cells.Item(11, 2).Value = 503
$ws4.Cells.Item(11, 3).Value = 201.12518310546875

$ws4.Cells.Item(12, 2).Value = 1985
$ws4.Cells.Item(12, 3).Value = 144.06742858886719

$ws4.Cells.Item(13, 2).Value = 14745
$ws4.Cells.Item(13, 3).Value = 620.298583984375
$ws4.Cells.Item(13, 4).Value = 826
$ws4.Cells.Item(13, 5).Value = 2138

$ws4.Cells.Item(14, 2).Value = 1658
$ws4.Cells.Item(14, 3).Value = 557.23602294921875

$ws4.Cells.Item(15, 2).Value = 2221
$ws4.Cells.Item(15, 3).Value = 578.8179931640625

$ws4.Cells.Item(16, 2).Value = 633
$ws4.Cells.Item(16, 3).Value = 224.13902282714844

$ws4.Cells.Item(17, 2).Value = 532
$ws4.Cells.Item(17, 3).Value = 195.77825927734375

$ws4.Cells.Item(18, 2).Value = 860
$ws4.Cells.Item(18, 3).Value = 350.52395629882813
$ws4.Cells.Item(18, 5).Value = 89

$ws4.Cells.Item(20, 2).Value = 7973
$ws4.Cells.Item(20, 3).Value = 461.96694946289063
$ws4.Cells.Item(20, 4).Value = 352

$ws4.Cells.Item(21, 2).Value = 1957
$ws4.Cells.Item(21, 3).Value = 642.0498046875

$ws4.Cells.Item(22, 2).Value = 2253
$ws4.Cells.Item(22, 3).Value = 484.00091552734375

# New narrower-than-default column C.
$ws4.Columns.Item(3).ColumnWidth = 15.592447916666666

# ---------------------------------------------------------------------
# Sheet 5: "Totalt antal per kön"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Cells.Item(2, 2).Value = 17821
$ws5.Cells.Item(2, 3).Value = 1603
$ws5.Cells.Item(2, 4).Value = 2556

$ws5.Cells.Item(3, 2).Value = 26908
$ws5.Cells.Item(3, 3).Value = 570
$ws5.Cells.Item(3, 4).Value = 2103

# ---------------------------------------------------------------------
# Sheet 6: "Totalt antal per åldersgrupp"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Cells.Item(2, 2).Value = 215
$ws6.Cells.Item(3, 2).Value = 608
$ws6.Cells.Item(4, 2).Value = 4479
$ws6.Cells.Item(5, 2).Value = 5963
$ws6.Cells.Item(5, 3).Value = 97
$ws6.Cells.Item(6, 2).Value = 6994
$ws6.Cells.Item(6, 3).Value = 248
$ws6.Cells.Item(7, 2).Value = 8313
$ws6.Cells.Item(7, 3).Value = 571
$ws6.Cells.Item(8, 2).Value = 5308
$ws6.Cells.Item(8, 3).Value = 652
$ws6.Cells.Item(9, 2).Value = 4431
$ws6.Cells.Item(9, 3).Value = 423
$ws6.Cells.Item(9, 4).Value = 1023
$ws6.Cells.Item(10, 2).Value = 5430
$ws6.Cells.Item(10, 4).Value = 1916
$ws6.Cells.Item(11, 2).Value = 2975

# ---------------------------------------------------------------------
# Sheet 7: rename the daily FOHM snapshot tab
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)
$ws7.Name = "FOHM  7 Jun 2020"

# ---------------------------------------------------------------------
# Selections / scroll anchors on the non-active sheets. Each of these
# sheets must be activated momentarily to record its own Selection, so
# we do them before re-activating sheet 2 (the sheet that stays active).
# ---------------------------------------------------------------------
$ws3.Activate()
$ws3.Columns.Item(2).Select()

$ws4.Activate()
$ws4.Columns.Item(5).Select()

$ws5.Activate()
$ws5.Columns.Item(4).Select()

$ws6.Activate()
$ws6.Range("F13").Select()

# Sheet 2 remains the active tab (workbook activeTab = 1), with its
# selection moved from B89 down to A91.
$ws2.Activate()
$ws2.Range("A91").Select()
